# Upgrade left table: add 2023 column (K) with data, matching existing
# table formatting, and tidy up the header row's border.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Bring the new column's formatting in line with the rest of the
#        table by copying the last existing data column (J) formats over,
#        then writing this municipality's 2023 figures.
$ws.Range("J3:J6").Copy() | Out-Null
$ws.Range("K3:K6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 4513
$ws.Range("K5").Value = 443
$ws.Range("K6").Value = 4070

# --- 2. Close the table off on the right with a thin border for the new
#        column, matching the top/bottom rules already used by the table.
$ws.Range("K3:K6").Borders.Item(10).LineStyle = 1
$ws.Range("K3:K6").Borders.Item(10).Weight = 2

$ws.Range("K3").Borders.Item(8).LineStyle = 1
$ws.Range("K3").Borders.Item(8).Weight = 2
$ws.Range("K3").Borders.Item(9).LineStyle = -4142

$ws.Range("K4").Borders.Item(8).LineStyle = 1
$ws.Range("K4").Borders.Item(8).Weight = 2

$ws.Range("K5").Borders.Item(8).LineStyle = -4142
$ws.Range("K5").Borders.Item(9).LineStyle = -4142

$ws.Range("K6").Borders.Item(8).LineStyle = -4142
$ws.Range("K6").Borders.Item(9).LineStyle = 1
$ws.Range("K6").Borders.Item(9).Weight = 2

# --- 3. The header row (B3:J3) no longer boxes itself in with a bottom
#        border now that the table continues below it.
$ws.Range("B3:J3").Borders.Item(9).LineStyle = -4142

# --- 4. Slightly widen the columns that sit under the new data so the
#        wider "2023" figures keep lining up like the rest of the sheet.
$ws.Range("G1:N1").ColumnWidth = 7.83
